$wb = $excel.ActiveWorkbook

# Map of sheet name -> cell address -> (old fragment, new fragment)
# We replace the "id='...'" attribute with "class='...'" in the
# "!!ObjTables type='Data' id='...'" header cells.

$targets = @(
    @{ Sheet = "!!FirstUnambiguousModel";  Cell = "A2" },
    @{ Sheet = "!!SecondUnambiguousModel"; Cell = "A1" },
    @{ Sheet = "!!TestModel";              Cell = "A1" },
    @{ Sheet = "!!TestModels3";            Cell = "A1" }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    $cell = $ws.Range($t.Cell)
    $text = $cell.Text
    if ($text -ne $null -and $text.Contains("!!ObjTables") -and $text.Contains("id='")) {
        $cell.Value = $text.Replace("id='", "class='")
    }
}
